$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 53; this shifts the existing rows 53-112
# down to 54-113 (and the sheet dimension grows to A1:R113 automatically).
$ws.Rows.Item(53).Insert()

# Populate the newly inserted (blank) row 53 with the new weekly record.
$ws.Range("A53").Value = 4
$ws.Range("B53").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C53").Value = 'Los Lagos'
$ws.Range("D53").Value = 44904
$ws.Range("E53").Value = 10
$ws.Range("F53").Value = 100112031
$ws.Range("G53").Value = 'Poroto verde'
$ws.Range("H53").Value = 'Magnum'
$ws.Range("I53").Value = 'Primera'
$ws.Range("J53").Value = 35
$ws.Range("K53").Value = 50000
$ws.Range("L53").Value = 50000
$ws.Range("M53").Value = 50000
$ws.Range("N53").Value = '$/malla 25 kilos'
$ws.Range("O53").Value = 'Provincia de Limarí'
$ws.Range("P53").Value = 2000
$ws.Range("Q53").Value = 25
$ws.Range("R53").Value = 'Hortaliza'
